$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value and force it to be stored as text, even if it
# looks like a pure number (e.g. "1.00" or "537.03"), matching the source
# data which stores every Price/Volume cell as a text string. Restoring the
# style back to "Normal" afterwards avoids leaving a stray number format on
# the cell once the text is committed.
function Set-TextValue($range, $val) {
    $range.NumberFormat = "@"
    $range.Value = $val
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '70.456.81'
$ws.Range('E2').Value = '  -2.79%  '

# Row 3
Set-TextValue $ws.Range('D3') '3.916.73'
$ws.Range('E3').Value = '  -3.04%  '

# Row 4
Set-TextValue $ws.Range('D4') '1.00'
$ws.Range('E4').Value = '  -0.07%  '

# Row 5
Set-TextValue $ws.Range('D5') '537.03'
$ws.Range('E5').Value = '  +3.94%  '

# Row 6
Set-TextValue $ws.Range('D6') '146.57'
$ws.Range('E6').Value = '  -0.14%  '

# Row 7
Set-TextValue $ws.Range('D7') '0.682'
$ws.Range('E7').Value = '  -5.75%  '

# Row 8
$ws.Range('E8').Value = '  +0.07%  '

# Row 9
Set-TextValue $ws.Range('D9') '0.730'
$ws.Range('E9').Value = '  -3.91%  '

# Row 10
Set-TextValue $ws.Range('D10') '0.166'
$ws.Range('E10').Value = '  -4.44%  '

# Row 11
Set-TextValue $ws.Range('D11') '52.94'
$ws.Range('E11').Value = '  +12.83%  '

# Row 12
Set-TextValue $ws.Range('D12') '0.0000314'
$ws.Range('E12').Value = '  -3.37%  '

# Row 13
Set-TextValue $ws.Range('D13') '10.43'
$ws.Range('E13').Value = '  -4.88%  '

# Row 14
Set-TextValue $ws.Range('D14') '4.548.53'
$ws.Range('E14').Value = '  -3.03%  '

# Row 15
Set-TextValue $ws.Range('D15') '3.927.84'
$ws.Range('E15').Value = '  -2.82%  '

# Row 16
Set-TextValue $ws.Range('D16') '13.79'
$ws.Range('E16').Value = '  -2.66%  '

# Row 17
Set-TextValue $ws.Range('D17') '20.17'
$ws.Range('E17').Value = '  -4.41%  '

# Row 18
Set-TextValue $ws.Range('D18') '0.131'
$ws.Range('E18').Value = '  -0.88%  '

# Row 19
Set-TextValue $ws.Range('D19') '1.16'
$ws.Range('E19').Value = '  -4.34%  '

# Row 20
Set-TextValue $ws.Range('D20') '70.420.68'
$ws.Range('E20').Value = '  -2.68%  '

# Row 21
Set-TextValue $ws.Range('D21') '426.06'
$ws.Range('E21').Value = '  -3.68%  '

# Row 22
Set-TextValue $ws.Range('D22') '95.92'
$ws.Range('E22').Value = '  -7.80%  '

# Row 23
Set-TextValue $ws.Range('D23') '3.48'
$ws.Range('E23').Value = '  -3.32%  '

# Row 24
Set-TextValue $ws.Range('D24') '4.16'
$ws.Range('E24').Value = '  +4.83%  '

# Row 25
Set-TextValue $ws.Range('D25') '14.07'
$ws.Range('E25').Value = '  -3.44%  '

# Row 26
Set-TextValue $ws.Range('D26') '11.10'
$ws.Range('E26').Value = '  -3.82%  '

# Row 27
Set-TextValue $ws.Range('D27') '10.45'
$ws.Range('E27').Value = '  -5.67%  '

# Row 28
$ws.Range('E28').Value = '  +0.82%  '

# Row 29
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D29') '36.07'
$ws.Range('E29').Value = '  -5.25%  '

# Row 30
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range('D30') '3.57'
$ws.Range('E30').Value = '  +14.86%  '

# Row 31
Set-TextValue $ws.Range('D31') '7.33'
$ws.Range('E31').Value = '  +7.76%  '

# Row 32
Set-TextValue $ws.Range('D32') '13.26'
$ws.Range('E32').Value = '  -3.00%  '

# Row 33
Set-TextValue $ws.Range('D33') '676.20'
$ws.Range('E33').Value = '  +1.13%  '

# Row 34
Set-TextValue $ws.Range('D34') '0.127'
$ws.Range('E34').Value = '  -0.34%  '

# Row 35
Set-TextValue $ws.Range('D35') '47.18'
$ws.Range('E35').Value = '  +11.91%  '

# Row 36
Set-TextValue $ws.Range('D36') '64.61'
$ws.Range('E36').Value = '  -4.60%  '

# Row 37
Set-TextValue $ws.Range('D37') '0.426'
$ws.Range('E37').Value = '  -1.76%  '

# Row 38
$ws.Range('E38').Value = '  -5.30%  '

# Row 39
Set-TextValue $ws.Range('D39') '3.41'
$ws.Range('E39').Value = '  -3.58%  '

# Row 40
Set-TextValue $ws.Range('D40') '0.146'
$ws.Range('E40').Value = '  -2.98%  '

# Row 41
$ws.Range('E41').Value = '  +0.20%  '

# Row 42
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range('D42') '0.999'
$ws.Range('E42').Value = '  +0.13%  '

# Row 43
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range('D43') '3.27'
$ws.Range('E43').Value = '  +2.36%  '

# Row 44
Set-TextValue $ws.Range('D44') '0.0475'
$ws.Range('E44').Value = '  -3.57%  '

# Row 45
Set-TextValue $ws.Range('D45') '0.147'
$ws.Range('E45').Value = '  -6.66%  '

# Row 46
Set-TextValue $ws.Range('D46') '2.66'
$ws.Range('E46').Value = '  -2.98%  '

# Row 47
Set-TextValue $ws.Range('D47') '9.58'
$ws.Range('E47').Value = '  +5.35%  '

# Row 48
Set-TextValue $ws.Range('D48') '3.32'
$ws.Range('E48').Value = '  -4.84%  '

# Row 49
$ws.Range('B49').Value = 'FLOKI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
Set-TextValue $ws.Range('D49') '0.000273'
$ws.Range('E49').Value = '  +1.34%  '

# Row 50
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range('D50') '2.94'
$ws.Range('E50').Value = '  -4.18%  '

# Row 51
Set-TextValue $ws.Range('D51') '144.77'
$ws.Range('E51').Value = '  +1.36%  '
